$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 -- this shifts the existing "FAPs" row
# (old row 2) down to row 3 and the existing "MuSCs" row (old row 3) down
# to row 4, matching the new ordering of ECs / FAPs / MuSCs.
$ws.Rows("2:2").Insert()

# New row 2: "ECs" sending cluster row (brand new data, recalculated TPM)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gm13306"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.065096
$ws.Range("H2").Value = 0.195288
$ws.Range("I2").Value = 0.02794828919627058
$ws.Range("J2").Value = 0.02794828919627058
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.8043309063093335
$ws.Range("R2").Value = 7.238978156784
$ws.Range("S2").Value = 0.02794828919627058
$ws.Range("T2").Value = 0.02794828919627058

# Row 3: previously the "FAPs" row (now shifted down from row 2) -- update
# its recalculated specificity values now that a third sending cluster
# (ECs) participates in the normalisation.
$ws.Range("I3").Value = 0.6668047231933247
$ws.Range("J3").Value = 0.6668047231933247
$ws.Range("Q3").Value = 19.190142322164
$ws.Range("R3").Value = 172.711280899476
$ws.Range("S3").Value = 0.6668047231933247
$ws.Range("T3").Value = 0.6668047231933247

# Row 4: previously the "MuSCs" row (now shifted down from row 3) -- update
# its recalculated specificity values the same way.
$ws.Range("I4").Value = 0.3052469876104046
$ws.Range("J4").Value = 0.3052469876104047
$ws.Range("Q4").Value = 8.784780509056446
$ws.Range("R4").Value = 79.06302458150802
$ws.Range("S4").Value = 0.3052469876104046
$ws.Range("T4").Value = 0.3052469876104047
